$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation needs to be inserted for this market/category,
# ahead of the existing row 55. This shifts the existing rows 55-69 down
# to 56-70 (dimension grows from A1:R69 to A1:R70) and inserts the new
# row with formatting (e.g. the date number format on column D) copied
# from the row below, matching Excel's default row-insert behaviour.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly record. The
# market/category descriptive columns are constant across this sheet, so
# they are simply repeated; the observation-specific columns (date,
# volume, prices, origin region and $/Kg) carry the new values.
$ws.Cells.Item(55, 1).Value = 5
$ws.Cells.Item(55, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(55, 3).Value = "Maule"
$ws.Cells.Item(55, 4).Value = 44841
$ws.Cells.Item(55, 5).Value = 7
$ws.Cells.Item(55, 6).Value = 100112040
$ws.Cells.Item(55, 7).Value = "Cilantro"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 150
$ws.Cells.Item(55, 11).Value = 8000
$ws.Cells.Item(55, 12).Value = 8000
$ws.Cells.Item(55, 13).Value = 8000
$ws.Cells.Item(55, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(55, 15).Value = "Región del Maule"
$ws.Cells.Item(55, 16).Value = 222
$ws.Cells.Item(55, 17).Value = 36
$ws.Cells.Item(55, 18).Value = "Hortaliza"
